# Remove the row for id 1059 (Company_59) - this shifts all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(55).Delete()

# After the shift, the row that now sits at 70 (originally holding Company_79's data
# shifted from row 71) gets its name replaced with garbage/test data.
$ws.Range("B70").Value = "sdflkjsldkfnlksdfsdflkjsdflksdlfkjsldkfjn sdvlklsdjflsjdflsjdfl sdlkhjsdlifjsldjf lisjdfoijsdf olisjdfoijsdof"

# The row that now sits at 80 (shifted from row 81) loses its created_at value.
$ws.Range("F80").Clear()
